# CORE_holdings.xlsx update:
#  - bump the "as of" date in the confidential disclosure footnote
#  - refresh the Weight / Percent Change figures in the holdings table
#
# The sheet ships protected (sheetProtection), so unlock it for the
# duration of the edit and restore protection afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# --- Footnote text: 2021-05-07 -> 2021-05-10 -------------------------------
$ws.Range("A11").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-10 for illustrative purposes only and are subject to change."

# --- Holdings table: Weight (D) / Percent Change (E) refresh ---------------
$ws.Range("D2").Value = 0.498236613455935
$ws.Range("E2").Value = -0.000199150292087058

$ws.Range("D3").Value = 0.244799809069033
$ws.Range("E3").Value = -0.01948333092798393

$ws.Range("D4").Value = 0.09650150950097673
$ws.Range("E4").Value = -0.0180345969819653

$ws.Range("D5").Value = 0.1034710535687272
$ws.Range("E5").Value = -0.006992371957864174

$ws.Range("D6").Value = 0.0299504162050779
$ws.Range("E6").Value = -0.01618823529411761

$ws.Range("D7").Value = 0.02704059820025013
$ws.Range("E7").Value = -0.02499808883112908

$ws.Range("E8").Value = -0.008493421244264776

$ws.Protect()
